$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 128
$ws.Range("F6").Value = 570
$ws.Range("F7").Value = 50
$ws.Range("F12").Value = 729
$ws.Range("F13").Value = 769
$ws.Range("F15").Value = 4
$ws.Range("F16").Value = 1525
$ws.Range("F17").Value = 1525
$ws.Range("F18").Value = 895
$ws.Range("F19").Value = 31
$ws.Range("F22").Value = 343
$ws.Range("F26").Value = 6674
$ws.Range("F27").Value = 5047
$ws.Range("F28").Value = 5047
$ws.Range("F33").Value = 28
$ws.Range("F35").Value = 1299
$ws.Range("F36").Value = 197
$ws.Range("F37").Value = 254
$ws.Range("F38").Value = 620
$ws.Range("F39").Value = 19
$ws.Range("F41").Value = 259
$ws.Range("F42").Value = 150
$ws.Range("F43").Value = 151
$ws.Range("F44").Value = 65
$ws.Range("F46").Value = 99

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 20
$ws.Range("F9").Value = 2
$ws.Range("F18").Value = 247

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2471
$ws.Range("F4").Value = 207
$ws.Range("F5").Value = 69

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 207
$ws.Range("F8").Value = 69
$ws.Range("F10").Value = 570
$ws.Range("F11").Value = 50
$ws.Range("F16").Value = 729
$ws.Range("F17").Value = 769
$ws.Range("F18").Value = 1525
$ws.Range("F19").Value = 1525
$ws.Range("F20").Value = 895
$ws.Range("F21").Value = 31
$ws.Range("F24").Value = 343
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 6674
$ws.Range("F30").Value = 5047
$ws.Range("F31").Value = 5047
$ws.Range("F34").Value = 1299
$ws.Range("F35").Value = 197
$ws.Range("F36").Value = 2
$ws.Range("F37").Value = 254
$ws.Range("F39").Value = 620
$ws.Range("F41").Value = 19
$ws.Range("F45").Value = 259
$ws.Range("F46").Value = 151
$ws.Range("F47").Value = 65
$ws.Range("F48").Value = 99
$ws.Range("F50").Value = 247
